$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.880.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.276.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.68%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.635'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.52%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.663'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +17.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.87'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.74'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0970'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.105'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.60%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.618.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.89'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.887'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.268.48'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.791.00'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000101'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.60%  '
$ws.Range("E24").Value = '  +4.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.69%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -8.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.127'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.42%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '31.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +31.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0805'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.07%  '
$ws.Range("E37").Value = '  +4.15%  '
$ws.Range("E38").Value = '  +12.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.29%  '
$ws.Range("E40").Value = '  +3.46%  '
$ws.Range("E41").Value = '  +6.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.79%  '
$ws.Range("E44").Value = '  +9.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.37%  '
$ws.Range("E48").Value = '  +2.63%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("E51").Value = '  +4.64%  '
